$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.283.61"
$ws.Range("E2").Value = "  -1.94%  "

$ws.Range("D3").Value = "3.430.20"

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.18"
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.55"
$ws.Range("E6").Value = "  +3.46%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +1.16%  "

$ws.Range("E9").Value = "  +4.55%  "

$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("E11").Value = "  +2.85%  "

$ws.Range("D12").Value = "4.012.50"
$ws.Range("E12").Value = "  -1.52%  "

$ws.Range("E13").Value = "  +0.31%  "

$ws.Range("E14").Value = "  -3.43%  "

$ws.Range("D15").Value = "3.467.82"
$ws.Range("E15").Value = "  -0.07%  "

$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("D17").Value = "62.282.57"
$ws.Range("E17").Value = "  -1.85%  "

$ws.Range("E18").Value = "  +2.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.52"
$ws.Range("E19").Value = "  +0.73%  "

$ws.Range("E20").Value = "  -4.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.54"
$ws.Range("E21").Value = "  -1.86%  "

$ws.Range("E22").Value = "  +1.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.22"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").Value = "3.560.06"
$ws.Range("E25").Value = "  -1.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000113"
$ws.Range("E26").Value = "  -3.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.181"
$ws.Range("E27").Value = "  +0.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("E28").Value = "  +1.10%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.97"
$ws.Range("E30").Value = "  -3.34%  "

$ws.Range("E31").Value = "  -0.85%  "

$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("E33").Value = "  -1.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.26"
$ws.Range("E34").Value = "  -1.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.47"
$ws.Range("E35").Value = "  +2.95%  "

$ws.Range("E36").Value = "  +1.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.94"
$ws.Range("E37").Value = "  -2.63%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.96"
$ws.Range("E38").Value = "  -0.53%  "

$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.10"
$ws.Range("E39").Value = "  -0.97%  "

$ws.Range("D40").Value = "3.462.77"
$ws.Range("E40").Value = "  -1.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0784"
$ws.Range("E41").Value = "  +2.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.82"
$ws.Range("E42").Value = "  +1.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.778"
$ws.Range("E43").Value = "  -2.65%  "

$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.69"
$ws.Range("E45").Value = "  -2.75%  "

$ws.Range("E46").Value = "  -2.47%  "

$ws.Range("D47").Value = "2.540.92"
$ws.Range("E47").Value = "  -2.97%  "

$ws.Range("E48").Value = "  +1.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.22"
$ws.Range("E49").Value = "  -3.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.64"
$ws.Range("E50").Value = "  -2.31%  "

$ws.Range("E51").Value = "  +0.00%  "
